$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update "Наличие балкона/лоджии" column (I) values to the simplified
# "есть"/"нет" wording, replacing the old "балкон"/"лоджия" values.
$ws.Range("I2").Value = "есть"
$ws.Range("I3").Value = "нет"
$ws.Range("I4").Value = "есть"

# Move/restore the active selection to I2, matching the saved view state.
$ws.Range("I2").Select()
